$d = $word.ActiveDocument

# Locate the "Docente(s) Responsável(eis)" heading paragraph
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd("`r", "`n") -eq "Docente(s) Responsável(eis) ") {
        $target = $p
        break
    }
}

# Insert a new paragraph right after it, formatted as a bulleted list item
# containing the professor's name/id.
$target.Range.InsertParagraphAfter()
$newPara = $target.Next()
$newPara.Range.Text = "5701460 - Antonio Iacono"
$newPara.Style = "ListBullet"
